$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New helper column: "wrong_nums" header + zeroed counters for each data row
$ws.Range("C1").Value = "wrong_nums"
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 0

# Match the saved selection state
$ws.Range("D13").Select()
